$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8591253757476807
$ws.Range("B1").Value = 1.702115297317505
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.891545057296753
$ws.Range("E1").Value = 1.128319501876831
